$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (StdType) values, set row by row (this establishes the
#     order in which the new "StdType"/"Run"/"Calib"/"RunCalib" strings
#     are first introduced into the shared string table) ---
$ws.Range("E1").Value = "StdType"
$ws.Range("E2").Value = "Run"
$ws.Range("E3").Value = "Run"
$ws.Range("E4").Value = "Calib"
$ws.Range("E5").Value = "RunCalib"
$ws.Range("E6").Value = "Calib"
$ws.Range("E7").Value = "Calib"
$ws.Range("E8").Value = "Calib"
$ws.Range("E9").Value = "Calib"
$ws.Range("E10").Value = "Calib"
$ws.Range("E11").Value = "Calib"
$ws.Range("E12").Value = "Calib"
$ws.Range("E13").Value = "Calib"
$ws.Range("E14").Value = "Calib"
$ws.Range("E15").Value = "Calib"
$ws.Range("E16").Value = "Calib"
$ws.Range("E17").Value = "Calib"
$ws.Range("E18").Value = "Calib"
$ws.Range("E19").Value = "Calib"

# --- Column D (REGEX) updated values. Set in the same order the author
#     apparently used: rows 2-15 in order, then 18, 19, 17, 16. ---
$ws.Range("D2").Value = "UWC\\D*3"
$ws.Range("D3").Value = "UWQ\\D*1"
$ws.Range("D4").Value = "UW.*Arg\\D*7"
$ws.Range("D5").Value = "UW\\D*6220"
$ws.Range("D6").Value = "UW\\D*6250"
$ws.Range("D7").Value = "UW.*Ank\\D*10"
$ws.Range("D8").Value = "UW.*Ank\\D*11"
$ws.Range("D9").Value = "UW.*Ank\\D*7"
$ws.Range("D10").Value = "UW.*Ank\\D*8"
$ws.Range("D11").Value = "UW.*Ank\\D*4"
$ws.Range("D12").Value = "UW.*Ank\\D*9"
$ws.Range("D13").Value = "UW.*Ank\\D*1"
$ws.Range("D14").Value = "UW.*Ank\\D*2"
$ws.Range("D15").Value = "UW.*Ank\\D*3"
$ws.Range("D18").Value = "UW.*Ank\\D*12"
$ws.Range("D19").Value = "UW.*Ank\\D*5.*opq"
$ws.Range("D17").Value = "UW.*Ank\\D*6.*a"
$ws.Range("D16").Value = "UW.*Ank\\D*5.*cl"

# --- Column widths (closest achievable to target 256ths-based widths) ---
$ws.Columns.Item(1).ColumnWidth = 13.330729166666666
$ws.Columns.Item(2).ColumnWidth = 7.830729166666666
$ws.Columns.Item(3).ColumnWidth = 10.6640625
$ws.Columns.Item(4).ColumnWidth = 20.830729166666668

# --- Selection ---
[void]$ws.Range("L22").Select()
